# PRJ0018886_Hierarchy Viewer+ Time recordManager(PArtial changes)
#
# Update the "Project_Title" sheet so that:
#  - the header row (A1:B1) is bold
#  - the Engagement cell (A2) text is updated to reflect the new engagement name
#  - the sheet becomes the active/selected tab, with selection on B4
#  - columns auto-fit their new (wider) content

$wb = $excel.ActiveWorkbook

$wsUsers      = $wb.Worksheets.Item("Users")
$wsTimePeriod = $wb.Worksheets.Item("Time_Record_Period_Title")
$wsProject    = $wb.Worksheets.Item("Project_Title")
$wsUpdate     = $wb.Worksheets.Item("Update_Timer")

# Update the engagement name shown under the Project dropdown selection.
$wsProject.Range("A2").Value = "GE Healthcare-GE Healthcare Bio-Sciences AB-FVA-101397"

# Make the header row bold (matches the styling already used on other sheets).
$wsProject.Range("A1:B1").Font.Bold = $true

# Resize the columns to fit the new (longer) text (best-fit widths).
$wsProject.Columns.Item(1).ColumnWidth = 48.666666666666664
$wsProject.Columns.Item(2).ColumnWidth = 15

# Make Project_Title the active sheet (tabSelected) before moving the selection,
# so the selection change applies to the now-active sheet/window.
$wsProject.Activate()

# Move the selection / active cell.
$wsProject.Range("B4").Select()

# Keep the sheet print orientation explicit (portrait), as set on the sheet.
$wsProject.PageSetup.Orientation = 1
